$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns that change: B, C, D, F, G, H, L (in that order)
$cols = @("B", "C", "D", "F", "G", "H", "L")

# Each row: row number followed by the 7 new values for columns B, C, D, F, G, H, L
$data = @(
    @(2, 18.02622406214402, 10.8732378308616, 5.187993927289388, 28.57249813807272, 36.70424317113589, 15.82437198321952, 10.71811381589169),
    @(3, 17.44610520383242, 10.62600114821866, 5.20608871023122, 28.45934946125375, 36.43878964619378, 15.85407002280443, 10.69427126272141),
    @(4, 17.08364535133095, 10.46984292237912, 5.217710498784029, 28.40057899600273, 36.29172190607667, 15.87663060746644, 10.68192810282209),
    @(5, 16.93460504147481, 10.40516738212959, 5.2225756478598, 28.37933278750872, 36.23584472437645, 15.88690701212714, 10.67747850899123),
    @(6, 16.90978375814523, 10.39436693070215, 5.223391320348133, 28.37596839443767, 36.22681249939064, 15.88867864521276, 10.67677478395571),
    @(7, 17.08164041608306, 10.4689748203236, 5.2177755881498, 28.40028150485844, 36.29095185275151, 15.87676482147, 10.68186574077567),
    @(8, 17.82763779148594, 10.78892623682928, 5.194127172710992, 28.53127148033133, 36.60944350854182, 15.83371114608799, 10.70941822592683),
    @(9, 19.2312123944741, 11.3795622867559, 5.151786520788955, 28.87233276732546, 37.3576436945996, 15.78381334213869, 10.78151748113075),
    @(10, 20.21470569092348, 11.78832393635975, 5.123103555150152, 29.17294566916862, 37.97847973665603, 15.76846120073212, 10.84527015090297),
    @(11, 20.64969459318475, 11.96830052838838, 5.110574243439148, 29.32020247455027, 38.27529114354999, 15.76614965825021, 10.87655355656964),
    @(12, 20.81248118508529, 12.03555601433269, 5.105903774812715, 29.37743889670002, 38.3896558838163, 15.76594912805573, 10.88872199556906),
    @(13, 20.77751042011359, 12.02111185674264, 5.106906355271364, 29.36504710953168, 38.36493948935885, 15.76596226248858, 10.88608708078496),
    @(14, 20.66312671018996, 11.97385188195676, 5.110188518618786, 29.3248820458595, 38.28466113868996, 15.76611962147975, 10.8775482425174),
    @(15, 20.59280731680465, 11.9447857664632, 5.112208575643887, 29.30047050256941, 38.23574163826537, 15.76630396591334, 10.87235971567725),
    @(16, 20.18601526069758, 11.77643846946756, 5.123932771642817, 29.16353039901684, 37.95936350566053, 15.76870654557497, 10.8432711195018),
    @(17, 19.93317153319095, 11.67160461916023, 5.13125769563116, 29.08218764814626, 37.79343032507119, 15.77137937254051, 10.82600652434828),
    @(18, 19.78658255741576, 11.61074705148718, 5.135519649972809, 29.03639285046274, 37.69935080553302, 15.77335633988798, 10.81629157147592),
    @(19, 19.73675567826841, 11.59004680606286, 5.136971079632551, 29.02105887029216, 37.66773379805161, 15.7741011155843, 10.81303939042141),
    @(20, 19.96020837116441, 11.68282260378006, 5.130472892556766, 29.09074437577838, 37.81095406331087, 15.77104932296748, 10.82782214146526),
    @(21, 20.69677758755468, 11.98775794178167, 5.109222460137781, 29.33663982491079, 38.30818823394785, 15.766055066531, 10.88004761540128),
    @(22, 21.16683569917738, 12.18180226077912, 5.095765827513285, 29.50591538200412, 38.64458444133684, 15.76672528629686, 10.9160542464615),
    @(23, 20.91703909001526, 12.07872946281898, 5.102908536043866, 29.41479906054627, 38.46403261542311, 15.76600675421431, 10.89666745520245),
    @(24, 19.94798882307217, 11.67775277894075, 5.130827543864912, 29.08687285422834, 37.80302746926195, 15.77119716701373, 10.8270006442442),
    @(25, 18.85912482925665, 11.22401568030816, 5.16281253801179, 28.77116935401684, 37.14240429714759, 15.79358699148334, 10.76010040873632)
)

foreach ($entry in $data) {
    $row = $entry[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $value = $entry[$i + 1]
        $ws.Range("$col$row").Value = $value
    }
}
